$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")

# Zero out the benchmark dataset values on Sheet2
$ws2.Range("C2").Value = 0
$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 0
$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 0
$ws2.Range("E4").Value = 0

# Make Sheet2 the active sheet/tab
$ws2.Activate()

# Update selection on Sheet2 to E5
$ws2.Range("E5").Select()
